# Add the new "checkMSProjectIntegration" worksheet as the last tab,
# give it the two new labels (wired through sharedStrings automatically),
# and leave it as the active/selected sheet - mirroring the author's
# "Tested Successfully for MS Project Integration button in Dashboard" commit.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "checkMSProjectIntegration"

$newSheet.Range("A1").Value = "msProjectIntegrationTitle"
$newSheet.Range("A2").Value = "MS Project Integration"

$newSheet.Columns.Item(1).ColumnWidth = 27.25

$newSheet.Range("C4").Select()
